$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Cypher query text stored in cell B2 -----------------------
# (new version: WHERE clause split onto its own line, WITH/RETURN clauses
#  reshaped to project Age/Weight through a CASE-based integer coercion,
#  and an ORDER BY + LIMIT 100 tacked on at the end.)
$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Mixed Breed']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$ws.Range("B2").Value = $newQuery

# --- Column widths ----------------------------------------------------------
# Previously columns B:C shared one <col> entry at width 75.7109375.
# Now column B widens to 86 while column C keeps its original 75.7109375 -
# only touch column B so the shared <col min="2" max="3"> entry splits into
# two distinct entries without perturbing column C's stored width.
$ws.Columns.Item(2).ColumnWidth = 85.16666666666667

# --- Row height ---------------------------------------------------------
# Row 2 grows from 255 to 345 points to fit the longer query text.
$ws.Rows.Item(2).RowHeight = 345

# --- Scroll position ------------------------------------------------------
# The saved view no longer pins topLeftCell to A2; reset the window so the
# sheet opens scrolled to the top-left (A1) instead.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
